# Applies the "Updated cryptos list" data refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to keep a literal text value even when the text
    # looks like a plain number (e.g. "293.39"), without leaving the
    # cell on a custom/text number format afterwards.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "40.026.21"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.231.84"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "293.39"
$ws.Range("E5").Value = "  -0.87%  "
Set-TextValue $ws.Range("D6") "86.63"
$ws.Range("E6").Value = "  +5.67%  "
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.82%  "
Set-TextValue $ws.Range("D10") "30.70"
$ws.Range("E10").Value = "  +7.17%  "
Set-TextValue $ws.Range("D11") "0.0790"
$ws.Range("E11").Value = "  +2.29%  "
Set-TextValue $ws.Range("D12") "47.06"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("E14").Value = "  +3.80%  "
$ws.Range("D15").Value = "2.576.17"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "2.238.99"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").Value = "39.937.90"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D21") "11.04"
$ws.Range("E21").Value = "  +9.31%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.80"
$ws.Range("E22").Value = "  +2.23%  "
Set-TextValue $ws.Range("D23") "65.24"
$ws.Range("E23").Value = "  +1.10%  "
Set-TextValue $ws.Range("D24") "235.08"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +3.50%  "
Set-TextValue $ws.Range("D27") "1.86"
$ws.Range("E27").Value = "  +5.58%  "
Set-TextValue $ws.Range("D28") "22.77"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  +2.91%  "
Set-TextValue $ws.Range("D30") "9.29"
$ws.Range("E30").Value = "  +3.00%  "
Set-TextValue $ws.Range("D31") "33.24"
$ws.Range("E31").Value = "  +5.03%  "
Set-TextValue $ws.Range("D32") "152.62"
$ws.Range("E32").Value = "  +3.45%  "
Set-TextValue $ws.Range("D33") "0.999"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +2.56%  "
Set-TextValue $ws.Range("D35") "0.0719"
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("E36").Value = "  +3.20%  "
Set-TextValue $ws.Range("D37") "16.24"
$ws.Range("E37").Value = "  +10.21%  "
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("E39").Value = "  +2.03%  "
Set-TextValue $ws.Range("D40") "0.0998"
$ws.Range("E40").Value = "  +5.32%  "
$ws.Range("E41").Value = "  +6.19%  "
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").Value = "2.041.39"
$ws.Range("E43").Value = "  +7.24%  "
$ws.Range("E44").Value = "  +7.92%  "
$ws.Range("E45").Value = "  +5.45%  "
Set-TextValue $ws.Range("D46") "10.08"
$ws.Range("E46").Value = "  +12.04%  "
Set-TextValue $ws.Range("D47") "16.83"
$ws.Range("E47").Value = "  +5.02%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "2.463.91"
$ws.Range("E49").Value = "  +1.61%  "
Set-TextValue $ws.Range("D50") "71.00"
$ws.Range("E50").Value = "  +1.16%  "
Set-TextValue $ws.Range("D51") "89.17"
$ws.Range("E51").Value = "  +3.08%  "
